# Actualización automática 2025-07-10 17:25:09
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Widen column D slightly (13 -> 14)
$ws.Columns.Item(4).ColumnWidth = 13.17

# Row 4 - FREGADEROS DE COCINA
$ws.Range("D4").Value = 64.81999999999999
$ws.Range("E4").Value = 581.1800000000001
$ws.Range("F4").Value = 0.1003405572755418

# Row 6 - GRIFERIAS
$ws.Range("D6").Value = 40.74
$ws.Range("E6").Value = 66.07999999999998
$ws.Range("F6").Value = 0.381389252948886

# Row 14 - PANELES PVC
$ws.Range("D14").Value = 3576.01
$ws.Range("E14").Value = -3336.01
$ws.Range("F14").Value = 14.90004166666667

# Row 16 - PORCELANATO
$ws.Range("D16").Value = 11850.48
$ws.Range("E16").Value = 42870.75
$ws.Range("F16").Value = 0.2165609216020912

# Row 19 - TOTAL
$ws.Range("D19").Value = 23498.79
$ws.Range("E19").Value = 81714.08
$ws.Range("F19").Value = 0.2233452048214254
